$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column CL with header "14-nov" and the corresponding daily values,
# mirroring the structure of the existing "13-nov" column (CK).

$ws.Range("CL1").Value = "14-nov"

$ws.Range("CL2").Value = 9
$ws.Range("CL3").Value = 10
$ws.Range("CL4").Value = 7
$ws.Range("CL5").Value = 9
$ws.Range("CL6").Value = 11
$ws.Range("CL7").Value = 5
$ws.Range("CL8").Value = 15
$ws.Range("CL9").Value = 14
$ws.Range("CL10").Value = 15
$ws.Range("CL11").Value = 0

# Match styling from the adjacent CK column (header style vs. data style)
$ws.Range("CL1").NumberFormat = $ws.Range("CK1").NumberFormat
$ws.Range("CL2:CL11").NumberFormat = $ws.Range("CK2:CK11").NumberFormat
$ws.Range("CL2:CL11").HorizontalAlignment = $ws.Range("CK2:CK11").HorizontalAlignment

# Update the selection to mirror the new last cell, like Excel does after editing
$ws.Range("CL11").Select()
